# "add scaling relation typesetting"
#
# - rename the two existing sheets
# - insert a new "Selectivity" sheet between them
# - retype the scaling-relation labels (near/line/triangle/paral -> Near/Line/
#   Triangle/Parall., plus Island/Overlayer) on the (renamed) Formation Energy
#   sheet
# - populate the new Selectivity sheet with its header + data table
# - fix up the selections on each sheet and make Selectivity the active tab

$wb = $excel.ActiveWorkbook

$wsFormationEnergy = $wb.Worksheets.Item(1)   # was "formation_energy"
$wsBackup          = $wb.Worksheets.Item(2)   # was "Formation Energy"

# Rename the pre-existing sheet that already holds the "Formation Energy"
# name out of the way first, so the rename below doesn't collide with it.
$wsBackup.Name = "Formation Energy backup"
$wsFormationEnergy.Name = "Formation Energy"

# New sheet, placed right after "Formation Energy".
$wsSelectivity = $wb.Worksheets.Add($null, $wsFormationEnergy)
$wsSelectivity.Name = "Selectivity"

# Inserting the new sheet shifts the backup sheet's position, so re-resolve
# $wsBackup by name now rather than keep using the (now stale/positional)
# reference captured above.
$wsBackup = $wb.Worksheets.Item("Formation Energy backup")

# --- Selectivity sheet: header row ------------------------------------------
$wsSelectivity.Range("A1").Value = "ΔG(HOCO*)-ΔG(H*)"
$wsSelectivity.Range("B1").Value = "Near"
$wsSelectivity.Range("C1").Value = "Line"
$wsSelectivity.Range("F1").Value = "Island "
$wsSelectivity.Range("G1").Value = "Overlayer"
$wsSelectivity.Range("A2").Value = "Pd"
$wsSelectivity.Range("D1").Value = "Triangle"

# --- Formation Energy sheet: retype the header labels -----------------------
$wsFormationEnergy.Range("F1").Value = "Island"
$wsSelectivity.Range("E1").Value = "Parall."
$wsFormationEnergy.Range("B1").Value = "Near"
$wsFormationEnergy.Range("C1").Value = "Line"
$wsFormationEnergy.Range("D1").Value = "Triangle"
$wsFormationEnergy.Range("E1").Value = "Parall."
$wsFormationEnergy.Range("G1").Value = "Overlayer"

# --- Selectivity sheet: data rows -------------------------------------------
$selectivityRows = @(
    @("Sc", 0.298, -0.557, -0.669, -0.233, -0.66,  -0.624),
    @("Ti", 0.302, -0.053, -0.147, -0.295, -0.28,  -0.179),
    @("V",  0.32,  0.742,  0.741,  -0.039, 0.065,  0.073),
    @("Mn", 0.281, 0.56,   0.808,  0.183,  0.066,  -0.027),
    @("Fe", 0.289, 0.832,  0.806,  0.43,   0.04,   -1.277),
    @("Co", 0.302, 0.882,  0.97,   0.278,  0.473,  -0.604),
    @("Ni", 0.291, 1.079,  1.222,  0.279,  0.249,  -0.57),
    @("Cu", 0.326, 1.496,  1.688,  0.866,  0.946,  0.134),
    @("Zn", 0.336, 0.251,  0.251,  0.22,   0.099,  1.215),
    @("Y",  0.312, -0.322, -0.68,  -0.86,  -1.131, -6.894),
    @("Zr", 0.278, -0.277, -0.237, -0.185, -0.591, -0.45),
    @("Nb", 1.228, 0.305,  0.152,  -0.064, 0.073,  0.467),
    @("Mo", 0.338, 0.529,  0.623,  -0.032, 0.052,  0.549),
    @("Ru", 0.317, 0.624,  0.614,  0.65,   0.592,  0.423),
    @("Rh", 0.274, 0.623,  0.504,  0.566,  0.53,   0.394),
    @("Ag", 0.328, 1.789,  0.683,  1.09,   1.146,  0.038)
)

# Row 2 (Pd) header cell already set above; fill in its numeric values now.
$wsSelectivity.Cells.Item(2, 2).Value = 0.288
$wsSelectivity.Cells.Item(2, 3).Value = 0.288
$wsSelectivity.Cells.Item(2, 4).Value = 0.287
$wsSelectivity.Cells.Item(2, 5).Value = 0.288
$wsSelectivity.Cells.Item(2, 6).Value = 0.288
$wsSelectivity.Cells.Item(2, 7).Value = 0.288

$r = 3
foreach ($row in $selectivityRows) {
    $wsSelectivity.Cells.Item($r, 1).Value = $row[0]
    $wsSelectivity.Cells.Item($r, 2).Value = $row[1]
    $wsSelectivity.Cells.Item($r, 3).Value = $row[2]
    $wsSelectivity.Cells.Item($r, 4).Value = $row[3]
    $wsSelectivity.Cells.Item($r, 5).Value = $row[4]
    $wsSelectivity.Cells.Item($r, 6).Value = $row[5]
    $wsSelectivity.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# --- selections --------------------------------------------------------
$wsFormationEnergy.Range("J7").Select()
$wsBackup.Range("B1").Select()
$wsSelectivity.Range("H8").Select()

# Make "Selectivity" the active tab.
$wsSelectivity.Activate()
